$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has a table of "problemes rencontres" running down to row 33.
# Add a new entry in row 34 ("deux ckeditor") about Ckeditor HTML content
# needing the Twig |raw filter. Start by cloning row 33's look (styles +
# row height) so the new row matches the rest of the table exactly, then
# overwrite only the cells that actually change.
$ws.Range("A33:H33").Copy($ws.Range("A34:H34"))
$ws.Rows(34).RowHeight = 83.25

$ws.Range("B34").Value = "affichage du contenu de Ckeditor (html interpreté)"
$ws.Range("D34").Value = "fitre |raw"
$ws.Range("E34").Value = 42107

# Liens-sources cell: text + a live hyperlink to the Twig raw filter doc.
$ws.Hyperlinks.Add($ws.Range("H34"), "http://twig.sensiolabs.org/doc/filters/raw.html")
$ws.Range("H34").Value = "http://twig.sensiolabs.org/doc/filters/raw.html"
# Re-apply H33's formatting (the hyperlink insert can restyle the cell).
$ws.Range("H33").Copy()
$ws.Range("H34").PasteSpecial(-4122)

# Leave the sheet selection/view over the (now larger) used range.
$ws.Range("A5:H34").Select() | Out-Null
